$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Бригадиры")

# New foreman (бригадир) names to add to the top of the list (right after the header row).
$newNames = @(
    "Абдуллоев Бузургмехр Мамадамонович",
    "Dilo",
    "Озар",
    "Abd"
)

$insertCount = $newNames.Count

# Determine how many existing data rows (below header in row 1) currently hold names.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Shift the existing rows down by $insertCount, working from the bottom up so
# values are not overwritten before they are read. This avoids Range.Insert(),
# which would otherwise copy cell formatting/styles onto the shifted rows.
for ($r = $lastRow; $r -ge 2; $r--) {
    $val = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + $insertCount, 1).Value2 = $val
}

# Write the new names into the now-empty rows right after the header.
for ($i = 0; $i -lt $insertCount; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value2 = $newNames[$i]
}
